$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Row 2 (surrogate -> David Haden) ---
$ws.Range("B2").Value = "David"
$ws.Range("D2").Value = "Haden"
$ws.Range("F2").Value = "david@yopmail.com"

# --- Update existing Row 3 (surrogate -> Emma, new phone number) ---
$ws.Range("B3").Value = "Emma"
$ws.Range("D3").Value = "Haden"
$ws.Range("E3").Value = 9123477789
$ws.Range("F3").Value = "emma@yopmail.com"

# --- New Row 4 (Egg donor) ---
$ws.Range("A4").Value = "EGG_DONER"
$ws.Range("B4").Value = "Sophia"
$ws.Range("D4").Value = "Loren"
$ws.Range("E4").Value = 7777777777
$ws.Range("F4").Value = "shophia@yopmail.com"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "11/3/1988"

# --- New Row 5 (Sperm donor) ---
$ws.Range("A5").Value = "SPERM_DONER"
$ws.Range("B5").Value = "James"
$ws.Range("D5").Value = "william"
$ws.Range("E5").Value = 8888888888
$ws.Range("F5").Value = "james@yopmail.com"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "11/3/1988"

# --- Hyperlinks for the new email cells ---
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:shophia@yopmail.com", "", "", "shophia@yopmail.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:james@yopmail.com", "", "", "james@yopmail.com")

# --- Re-apply the big blue "email link" look to F4:F5 (sz 48, blue, Arial, no underline) ---
$linkRange = $ws.Range("F4:F5")
$linkRange.Font.Name = "Arial"
$linkRange.Font.Underline = $false
$linkRange.Font.Size = 48
$linkRange.Font.Color = 16711680

# --- Match the final selection shown in the source file ---
$null = $ws.Range("G5").Select()
